$d = $word.ActiveDocument

# 1. "GitHub" (4 runs: Git/H/u/b) -> single run "Github"
#    Use the Hyperlink object so the run keeps its Hyperlink character style
#    and merges into one run instead of four.
$hGithub = $d.Hyperlinks.Item(3)
$hGithub.TextToDisplay = "Github"

# 2. Merge "PyAudio" proofErr-wrapped run back into the surrounding text.
$pPyAudio = $d.Paragraphs.Item(22)
$rPyAudio = $pPyAudio.Range
$rPyAudio.Find.Execute("testing it with PyAudio, and recording", $true, $false, $false, $false, $false, `
                        $true, 1, $false, "testing it with PyAudio, and recording", 2) | Out-Null

# 3. Merge "aforementioned subjects" proofErr-wrapped run back into surrounding text.
$pTutor = $d.Paragraphs.Item(26)
$rTutor = $pTutor.Range
$rTutor.Find.Execute("the aforementioned subjects. All", $true, $false, $false, $false, $false, `
                      $true, 1, $false, "the aforementioned subjects. All", 2) | Out-Null

# 4. Merge "Machinery(" proofErr-wrapped run back into surrounding text.
$pACM = $d.Paragraphs.Item(32)
$rACM = $pACM.Range
$rACM.Find.Execute("Computing Machinery(ACM) chapter", $true, $false, $false, $false, $false, `
                    $true, 1, $false, "Computing Machinery(ACM) chapter", 2) | Out-Null

# 5. "Hackerrank (Java)" hyperlink: merge spellStart/Hackerrank/spellEnd/" (Java)" runs into one.
$hHackerrankJava = $d.Hyperlinks.Item(6)
$hHackerrankJava.TextToDisplay = "Hackerrank (Java)"

# 6. "Hackerrank (JavaScript)" hyperlink: same merge.
$hHackerrankJs = $d.Hyperlinks.Item(8)
$hHackerrankJs.TextToDisplay = "Hackerrank (JavaScript)"
